$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("products")

# Add SKU values to F3, G3, H3 (matching E3's value)
$ws.Range("F3").Value = 49313
$ws.Range("G3").Value = 49313
$ws.Range("H3").Value = 49313

# Update the active selection to H3
$ws.Range("H3").Select()
